# Apply the "stuff at the bottom of the sheets" commit.
#
# 1) Rows 2-5 gain a "pair_kind" value of "generic" in column J
#    (header row 1's J1 = "pair_kind" already existed; these data rows
#    simply didn't have a J value before).
# 2) A new "stim details" block is appended starting at row 27:
#      row 27: section title "stim details"
#      row 28: headers  month | word_type | need_audio | need_image | word | count | find images
#      rows 29-32: month=6,6,7,7  word_type=video
#      rows 33-36: month=6,6,7,7  word_type=audio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the missing pair_kind ("generic") for the practice rows ---
$ws.Cells.Item(2, 10).Value = "generic"
$ws.Cells.Item(3, 10).Value = "generic"
$ws.Cells.Item(4, 10).Value = "generic"
$ws.Cells.Item(5, 10).Value = "generic"

# --- new "stim details" section ---
$ws.Cells.Item(27, 1).Value = "stim details"

$ws.Cells.Item(28, 1).Value = "month"
$ws.Cells.Item(28, 2).Value = "word_type"
$ws.Cells.Item(28, 3).Value = "need_audio"
$ws.Cells.Item(28, 4).Value = "need_image"
$ws.Cells.Item(28, 5).Value = "word"
$ws.Cells.Item(28, 6).Value = "count"
$ws.Cells.Item(28, 7).Value = "find images"

$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "video"
$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "video"
$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "video"
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "video"

$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "audio"
$ws.Cells.Item(34, 1).Value = 6
$ws.Cells.Item(34, 2).Value = "audio"
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "audio"
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "audio"
